$wb = $excel.ActiveWorkbook
$wsContactos = $wb.Worksheets.Item("Contactos ")
$wsPonencia = $wb.Worksheets.Item("Ponencia")

# --- Text corrections: accents / capitalization fixes ---
$wsContactos.Range("A10").Value = "César Omar Rojas Ayala"
$wsContactos.Range("A12").Value = "Efraín Rivera Roldán"
$wsContactos.Range("A15").Value = "Isaías Chala Ibargüen"
$wsContactos.Range("A19").Value = "José Huber Araujo Nieto"
$wsContactos.Range("A21").Value = "Juan Pablo Gallo Maya"
$wsContactos.Range("A24").Value = "Marcos Pérez Jiménez"
$wsContactos.Range("A28").Value = "Rodrigo Armando Lara Sánchez"
$wsContactos.Range("G2").Value = "Atlántico/Barranquilla"
$wsContactos.Range("G3").Value = "Caquetá/Florencia"
$wsContactos.Range("G6").Value = "San Andrés/Alcaldía Municipal Providencia y Santa Catalina Islas"
$wsContactos.Range("G7").Value = "Guanía/Puerto Inírida"
$wsContactos.Range("G8").Value = "Quindío/Armenia "
$wsContactos.Range("G9").Value = "Cauca/Popayán "
$wsContactos.Range("G11").Value = "Vaupés/Mitú"
$wsContactos.Range("G15").Value = "Choco/Quibdó "
$wsContactos.Range("G20").Value = "Caldas/Manizales"
$wsContactos.Range("G22").Value = "Bolívar/Cartagena"
$wsContactos.Range("G23").Value = "Córdoba/Montería "
$wsContactos.Range("G25").Value = "Valle del cauca/Santiago de Cali"
$wsContactos.Range("G26").Value = "Boyacá/Tunja"

# --- Selection / active sheet changes ---
$wsPonencia.Activate()
$wsPonencia.Range("A3").Select()
$wsContactos.Activate()
$wsContactos.Range("A3").Select()
